# Applies the scheduled-runner market-data refresh to the Leve profit
# calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each edited row holds current market-board prices (H:currentAveragePrice,
# I:currentAveragePriceNQ, J:currentAveragePriceHQ, K:LevePriceNQ,
# L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) that get refreshed with
# the latest pull; a couple of rows lose a stale NQ/HQ profit cell entirely
# when that side no longer has data.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3050
$ws.Range("I69").Value = 2100
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 6300
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -5426
$ws.Range("N69").Value = -13748

$ws.Range("H72").Value = 3050
$ws.Range("I72").Value = 2100
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 18900
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -14532
$ws.Range("N72").Value = -44736

$ws.Range("H112").Value = 2029.9565
$ws.Range("J112").Value = 2286.45
$ws.Range("L112").Value = 6859.349999999999
$ws.Range("N112").Value = -9075.349999999999

$ws.Range("H132").Value = 2381889
$ws.Range("I132").Value = 2874454
$ws.Range("J132").Value = 1158
$ws.Range("K132").Value = 8623362
$ws.Range("L132").Value = 3474
$ws.Range("M132").Value = -8620832
$ws.Range("N132").Value = -8534

$ws.Range("H137").Value = 1620.8948
$ws.Range("I137").Value = 1458.4706
$ws.Range("J137").Value = 3001.5
$ws.Range("K137").Value = 4375.4118
$ws.Range("L137").Value = 9004.5
$ws.Range("M137").Value = -1825.4118
$ws.Range("N137").Value = -14104.5


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 6800
$ws.Range("J40").Value = 6800
$ws.Range("L40").Value = 6800
$ws.Range("N40").Value = -7152

$ws.Range("H61").Value = 3083.639
$ws.Range("I61").Value = 2939.32
$ws.Range("J61").Value = 3411.6365
$ws.Range("K61").Value = 2939.32
$ws.Range("L61").Value = 3411.6365
$ws.Range("M61").Value = -2727.32
$ws.Range("N61").Value = -3835.6365

$ws.Range("H74").Value = 886.9231
$ws.Range("I74").Value = 761.6667
$ws.Range("J74").Value = 1057.7273
$ws.Range("K74").Value = 761.6667
$ws.Range("L74").Value = 1057.7273
$ws.Range("M74").Value = 112.3333
$ws.Range("N74").Value = -2805.7273

$ws.Range("H77").Value = 886.9231
$ws.Range("I77").Value = 761.6667
$ws.Range("J77").Value = 1057.7273
$ws.Range("K77").Value = 3808.3335
$ws.Range("L77").Value = 5288.636500000001
$ws.Range("M77").Value = 559.6665000000003
$ws.Range("N77").Value = -14024.6365

$ws.Range("H132").Value = 32292064
$ws.Range("I132").Value = 47620230
$ws.Range("J132").Value = 102910.4
$ws.Range("K132").Value = 142860690
$ws.Range("L132").Value = 308731.2
$ws.Range("M132").Value = -142858160
$ws.Range("N132").Value = -313791.2

$ws.Range("H136").Value = 3083.639
$ws.Range("I136").Value = 2939.32
$ws.Range("J136").Value = 3411.6365
$ws.Range("K136").Value = 8817.960000000001
$ws.Range("L136").Value = 10234.9095
$ws.Range("M136").Value = -6267.960000000001
$ws.Range("N136").Value = -15334.9095


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4889.5933
$ws.Range("I134").Value = 2553.068
$ws.Range("J134").Value = 11743.4
$ws.Range("K134").Value = 7659.204000000001
$ws.Range("L134").Value = 35230.2
$ws.Range("M134").Value = -5124.204000000001
$ws.Range("N134").Value = -40300.2


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1606.6316
$ws.Range("I58").Value = 881.65216
$ws.Range("J58").Value = 2718.2666
$ws.Range("K58").Value = 881.65216
$ws.Range("L58").Value = 2718.2666
$ws.Range("M58").Value = -678.65216
$ws.Range("N58").Value = -3124.2666

$ws.Range("H86").Value = 16080.35
$ws.Range("I86").Value = 15000.077
$ws.Range("J86").Value = 18086.572
$ws.Range("K86").Value = 15000.077
$ws.Range("L86").Value = 18086.572
$ws.Range("M86").Value = -13877.077
$ws.Range("N86").Value = -20332.572

$ws.Range("H89").Value = 16080.35
$ws.Range("I89").Value = 15000.077
$ws.Range("J89").Value = 18086.572
$ws.Range("K89").Value = 75000.38499999999
$ws.Range("L89").Value = 90432.86
$ws.Range("M89").Value = -69384.38499999999
$ws.Range("N89").Value = -101664.86

$ws.Range("H132").Value = 4547620.5
$ws.Range("I132").Value = 1422.7142
$ws.Range("J132").Value = 22227278
$ws.Range("K132").Value = 4268.142599999999
$ws.Range("L132").Value = 66681834
$ws.Range("M132").Value = -1738.142599999999
$ws.Range("N132").Value = -66686894

$ws.Range("H134").Value = 6254.216
$ws.Range("I134").Value = 1427.8422
$ws.Range("J134").Value = 20362.076
$ws.Range("K134").Value = 4283.5266
$ws.Range("L134").Value = 61086.228
$ws.Range("M134").Value = -1748.5266
$ws.Range("N134").Value = -66156.228

$ws.Range("H136").Value = 1606.6316
$ws.Range("I136").Value = 881.65216
$ws.Range("J136").Value = 2718.2666
$ws.Range("K136").Value = 2644.95648
$ws.Range("L136").Value = 8154.7998
$ws.Range("M136").Value = -94.95647999999983
$ws.Range("N136").Value = -13254.7998


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 2001
$ws.Range("I13").Value = 2001
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 6003
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -5835
$ws.Range("N13").ClearContents()

$ws.Range("H132").Value = 1451.9445
$ws.Range("I132").Value = 858.3077
$ws.Range("J132").Value = 2995.4
$ws.Range("K132").Value = 7724.7693
$ws.Range("L132").Value = 26958.6
$ws.Range("M132").Value = -5194.7693
$ws.Range("N132").Value = -32018.6


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18467000
$ws.Range("I11").Value = 27575000
$ws.Range("J11").Value = 251000
$ws.Range("K11").Value = 27575000
$ws.Range("L11").Value = 251000
$ws.Range("M11").Value = -27574861
$ws.Range("N11").Value = -251278

$ws.Range("H132").Value = 44658.375
$ws.Range("I132").Value = 2590.0715
$ws.Range("K132").Value = 7770.2145
$ws.Range("M132").Value = -5240.2145


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H132").Value = 2979647.8
$ws.Range("I132").Value = 4632652
$ws.Range("J132").Value = 4240.3
$ws.Range("K132").Value = 13897956
$ws.Range("L132").Value = 12720.9
$ws.Range("M132").Value = -13895426
$ws.Range("N132").Value = -17780.9

$ws.Range("H136").Value = 22226558
$ws.Range("I136").Value = 3997.1428
$ws.Range("K136").Value = 11991.4284
$ws.Range("M136").Value = -9441.428400000001


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4512.857
$ws.Range("I81").Value = 4296.6665
$ws.Range("J81").Value = 4675
$ws.Range("K81").Value = 8593.333000000001
$ws.Range("L81").Value = 9350
$ws.Range("M81").Value = -7532.333000000001
$ws.Range("N81").Value = -11472

$ws.Range("H84").Value = 4512.857
$ws.Range("I84").Value = 4296.6665
$ws.Range("J84").Value = 4675
$ws.Range("K84").Value = 42966.665
$ws.Range("L84").Value = 46750
$ws.Range("M84").Value = -37662.665
$ws.Range("N84").Value = -57358

$ws.Range("H132").Value = 51852948
$ws.Range("I132").Value = 70312920
$ws.Range("J132").Value = 6413001.5
$ws.Range("K132").Value = 210938760
$ws.Range("L132").Value = 19239004.5
$ws.Range("M132").Value = -210936230
$ws.Range("N132").Value = -19244064.5

$ws.Range("H136").Value = 42492.543
$ws.Range("I136").Value = 50801.25
$ws.Range("K136").Value = 152403.75
$ws.Range("M136").Value = -149853.75

